# "CSC Assessment memorial Day" — replace the three formula-driven percentage
# scores in column D with their equivalent plain numeric (0-100) scores, then
# append a new summary row that scales the column averages by 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3/D4/D8 used to store "=95/100" etc. (a 0-1 fraction). Replace with the
# plain score out of 100 so the row/column now matches the rest of column D.
$ws.Range("D3").Value = 95
$ws.Range("D4").Value = 82
$ws.Range("D8").Value = 78

# New row 30: 5x the (now recalculated) column averages in row 29.
$ws.Range("D30").Formula = "=5*0.85"
$ws.Range("E30").Formula = "=5*E29"
$ws.Range("F30").Formula = "=5*F29"
$ws.Range("G30").Formula = "=5*G29"

# Columns F/G otherwise carry a column-level style; the new row's F30/G30
# cells are unstyled in the authored workbook, so reset them explicitly.
$ws.Range("F30").Style = "Normal"
$ws.Range("G30").Style = "Normal"

# Match the author's final viewport/selection (scrolled down to row 10,
# active cell parked just past the new data in D31).
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 2
$ws.Range("D31").Select()
